$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -5972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4374.1577
$ws.Range("I32").Value = 4472.108
$ws.Range("K32").Value = 4472.108
$ws.Range("M32").Value = -4185.108

$ws.Range("H61").Value = 1326.875
$ws.Range("I61").Value = 769.1667
$ws.Range("K61").Value = 769.1667
$ws.Range("M61").Value = -557.1667

$ws.Range("H74").Value = 643.9655
$ws.Range("I74").Value = 617.1111
$ws.Range("K74").Value = 617.1111
$ws.Range("M74").Value = 256.8889

$ws.Range("H76").Value = 13999.5
$ws.Range("J76").Value = 13999.5
$ws.Range("L76").Value = 13999.5
$ws.Range("N76").Value = -14675.5

$ws.Range("H77").Value = 643.9655
$ws.Range("I77").Value = 617.1111
$ws.Range("K77").Value = 3085.5555
$ws.Range("M77").Value = 1282.4445

$ws.Range("H79").Value = 13999.5
$ws.Range("J79").Value = 13999.5
$ws.Range("L79").Value = 13999.5
$ws.Range("N79").Value = -16339.5

$ws.Range("H136").Value = 1326.875
$ws.Range("I136").Value = 769.1667
$ws.Range("K136").Value = 2307.5001
$ws.Range("M136").Value = 242.4998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 791.8461
$ws.Range("I80").Value = 562
$ws.Range("J80").Value = 860.8
$ws.Range("K80").Value = 562
$ws.Range("L80").Value = 860.8
$ws.Range("M80").Value = 436
$ws.Range("N80").Value = -2856.8

$ws.Range("H83").Value = 791.8461
$ws.Range("I83").Value = 562
$ws.Range("J83").Value = 860.8
$ws.Range("K83").Value = 2810
$ws.Range("L83").Value = 4304
$ws.Range("M83").Value = 2182
$ws.Range("N83").Value = -14288

$ws.Range("H103").Value = 5657
$ws.Range("J103").Value = 5657
$ws.Range("L103").Value = 5657
$ws.Range("N103").Value = -8001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 62501108
$ws.Range("I16").Value = 100000960
$ws.Range("J16").Value = 1365.1666
$ws.Range("K16").Value = 100000960
$ws.Range("L16").Value = 1365.1666
$ws.Range("M16").Value = -100000673
$ws.Range("N16").Value = -1939.1666

$ws.Range("H31").Value = 1232.0541
$ws.Range("I31").Value = 862.9524
$ws.Range("K31").Value = 862.9524
$ws.Range("M31").Value = -567.9524

$ws.Range("H34").Value = 1232.0541
$ws.Range("I34").Value = 862.9524
$ws.Range("K34").Value = 862.9524
$ws.Range("M34").Value = -660.9524

$ws.Range("H58").Value = 797.0526
$ws.Range("I58").Value = 750.26666
$ws.Range("K58").Value = 750.26666
$ws.Range("M58").Value = -547.26666

$ws.Range("H113").Value = 62501108
$ws.Range("I113").Value = 100000960
$ws.Range("J113").Value = 1365.1666
$ws.Range("K113").Value = 100000960
$ws.Range("L113").Value = 1365.1666
$ws.Range("M113").Value = -99998790
$ws.Range("N113").Value = -5705.1666

$ws.Range("H122").Value = 1200
$ws.Range("J122").Value = 1200
$ws.Range("L122").Value = 3600
$ws.Range("N122").Value = -8500

$ws.Range("H134").Value = 9260639
$ws.Range("I134").Value = 11905968
$ws.Range("K134").Value = 35717904
$ws.Range("M134").Value = -35715369

$ws.Range("H136").Value = 797.0526
$ws.Range("I136").Value = 750.26666
$ws.Range("K136").Value = 2250.79998
$ws.Range("M136").Value = 299.2000200000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 450
$ws.Range("J29").Value = 450
$ws.Range("L29").Value = 1350
$ws.Range("N29").Value = -1904

$ws.Range("H50").Value = 283.22223
$ws.Range("I50").Value = 237.25
$ws.Range("J50").Value = 320
$ws.Range("K50").Value = 711.75
$ws.Range("L50").Value = 960
$ws.Range("M50").Value = -230.75
$ws.Range("N50").Value = -1922

$ws.Range("H53").Value = 283.22223
$ws.Range("I53").Value = 237.25
$ws.Range("J53").Value = 320
$ws.Range("K53").Value = 711.75
$ws.Range("L53").Value = 960
$ws.Range("M53").Value = -230.75
$ws.Range("N53").Value = -1922

$ws.Range("H70").Value = 3562.5
$ws.Range("J70").Value = 5400
$ws.Range("L70").Value = 16200
$ws.Range("N70").Value = -16830

$ws.Range("H73").Value = 3562.5
$ws.Range("J73").Value = 5400
$ws.Range("L73").Value = 16200
$ws.Range("N73").Value = -18384

$ws.Range("H76").Value = 6411.8
$ws.Range("J76").Value = 6586.4287
$ws.Range("L76").Value = 19759.2861
$ws.Range("N76").Value = -20525.2861

$ws.Range("H79").Value = 6411.8
$ws.Range("J79").Value = 6586.4287
$ws.Range("L79").Value = 19759.2861
$ws.Range("N79").Value = -22411.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3729.9375
$ws.Range("I122").Value = 3544.75
$ws.Range("K122").Value = 10634.25
$ws.Range("M122").Value = -8184.25

$ws.Range("H126").Value = 2082.9375
$ws.Range("I126").Value = 1618.9166
$ws.Range("K126").Value = 4856.7498
$ws.Range("M126").Value = -2386.7498

$ws.Range("H132").Value = 2850.318
$ws.Range("I132").Value = 2710.0908
$ws.Range("J132").Value = 2990.5454
$ws.Range("K132").Value = 8130.2724
$ws.Range("L132").Value = 8971.636200000001
$ws.Range("M132").Value = -5600.2724
$ws.Range("N132").Value = -14031.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 569.94116
$ws.Range("I55").Value = 98.71429000000001
$ws.Range("J55").Value = 899.8
$ws.Range("K55").Value = 98.71429000000001
$ws.Range("L55").Value = 899.8
$ws.Range("M55").Value = 74.28570999999999
$ws.Range("N55").Value = -1245.8

$ws.Range("H132").Value = 61636
$ws.Range("J132").Value = 85976.086
$ws.Range("L132").Value = 257928.258
$ws.Range("N132").Value = -262988.258

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 4400
$ws.Range("J63").Value = 5800
$ws.Range("L63").Value = 5800
$ws.Range("N63").Value = -7048

$ws.Range("H66").Value = 4400
$ws.Range("J66").Value = 5800
$ws.Range("L66").Value = 17400
$ws.Range("N66").Value = -23640

$ws.Range("H69").Value = 28999.666
$ws.Range("J69").Value = 28999.666
$ws.Range("L69").Value = 28999.666
$ws.Range("N69").Value = -30497.666

$ws.Range("H72").Value = 28999.666
$ws.Range("J72").Value = 28999.666
$ws.Range("L72").Value = 86998.99800000001
$ws.Range("N72").Value = -94486.99800000001

$ws.Range("H95").Value = 29999.5
$ws.Range("J95").Value = 29999.5
$ws.Range("L95").Value = 29999.5
$ws.Range("N95").Value = -35491.5

$ws.Range("H132").Value = 2494.65
$ws.Range("I132").Value = 2043.125
$ws.Range("J132").Value = 4300.75
$ws.Range("K132").Value = 6129.375
$ws.Range("L132").Value = 12902.25
$ws.Range("M132").Value = -3599.375
$ws.Range("N132").Value = -17962.25

$ws.Range("H136").Value = 733.3333
$ws.Range("I136").Value = 527.36365
$ws.Range("K136").Value = 1582.09095
$ws.Range("M136").Value = 967.90905
